# Auto-generated Excel COM-interop script
# Applies updated crypto price/volume data per the target diff.
# All D-column price cells are plain text in the source workbook (many
# contain thousand-separator dots like '93.429.78'), so every D-cell write
# is forced to Text format first to avoid Excel reinterpreting the string
# as a number (which would introduce float rounding / drop trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '93.001.55'
$ws.Range('E2').Value = '  -5.06%  '

$ws.Range('D3').Value = '3.417.49'
$ws.Range('E3').Value = '  +1.60%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '234.34'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -7.38%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '635.85'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -3.71%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.41'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.26%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.392'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -8.78%  '

$ws.Range('E9').Value = '  +0.13%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.946'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -7.58%  '

$ws.Range('D11').Value = '3.412.95'
$ws.Range('E11').Value = '  +1.50%  '

$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.196'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -6.10%  '

$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '41.35'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.76%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.08'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.98%  '

$ws.Range('D15').Value = '93.038.08'
$ws.Range('E15').Value = '  -4.81%  '

$ws.Range('D16').Value = '4.059.92'
$ws.Range('E16').Value = '  +1.73%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000247'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -3.87%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '8.23'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -6.75%  '

$ws.Range('D19').Value = '3.419.19'
$ws.Range('E19').Value = '  +1.78%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.39'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -3.15%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.13'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.08%  '

$ws.Range('E22').Value = '  -9.92%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '490.91'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -4.64%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.18'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -5.91%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.0000190'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -5.85%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.33'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -8.13%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '90.36'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -6.95%  '

$ws.Range('D28').Value = '3.602.35'
$ws.Range('E28').Value = '  +1.73%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '11.80'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.91%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '11.53'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.00%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.04%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.71'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +4.29%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.134'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -8.70%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.177'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -7.55%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.11%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '29.74'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +3.14%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.546'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -4.61%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '546.72'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +3.96%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.43'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -6.58%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '7.49'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -6.87%  '

$ws.Range('E41').Value = '  -0.05%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.150'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.95%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.904'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +4.86%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '23.98'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.89%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.69'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -3.71%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0406'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -9.99%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '5.46'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.54%  '

$ws.Range('E48').Value = '  -2.89%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.12'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.06%  '

$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '52.78'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.36%  '

$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '3.17'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.02%  '

